$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)
$blankStyle = $ws.Cells.Item(2, 1).Style

function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $blankStyle
}

Set-TextCell $ws.Cells.Item(31, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(31, 2) "80001841"
Set-TextCell $ws.Cells.Item(31, 3) "3012"
Set-TextCell $ws.Cells.Item(31, 4) "Гострини на розрізі контакту"
Set-TextCell $ws.Cells.Item(31, 5) "1201"

Set-TextCell $ws.Cells.Item(32, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(32, 2) "80001841"
Set-TextCell $ws.Cells.Item(32, 3) "3070"
Set-TextCell $ws.Cells.Item(32, 4) "Гострини на розрізі контакту"
Set-TextCell $ws.Cells.Item(32, 5) "1202"

Set-TextCell $ws.Cells.Item(33, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(33, 2) "80001841"
Set-TextCell $ws.Cells.Item(33, 3) "3070"
Set-TextCell $ws.Cells.Item(33, 4) "Гострини на розрізі контакту"
Set-TextCell $ws.Cells.Item(33, 5) "1203"

Set-TextCell $ws.Cells.Item(34, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(34, 2) "80001841"
Set-TextCell $ws.Cells.Item(34, 3) "3070"
Set-TextCell $ws.Cells.Item(34, 4) "Гострини на розрізі контакту"
Set-TextCell $ws.Cells.Item(34, 5) "1204"

Set-TextCell $ws.Cells.Item(35, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(35, 2) "80001841"
Set-TextCell $ws.Cells.Item(35, 3) "3070"
Set-TextCell $ws.Cells.Item(35, 4) "Гострини на розрізі контакту"
Set-TextCell $ws.Cells.Item(35, 5) "1204"

Set-TextCell $ws.Cells.Item(36, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(36, 2) "80001841"
Set-TextCell $ws.Cells.Item(36, 3) "3070"
Set-TextCell $ws.Cells.Item(36, 4) "Гострини на розрізі контакту"
Set-TextCell $ws.Cells.Item(36, 5) "1205"

Set-TextCell $ws.Cells.Item(37, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(37, 2) "80001841"
Set-TextCell $ws.Cells.Item(37, 3) "3070"
Set-TextCell $ws.Cells.Item(37, 4) "Пошкодження поверхні контакту"
Set-TextCell $ws.Cells.Item(37, 5) "1206"

Set-TextCell $ws.Cells.Item(38, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(38, 2) "80001841"
Set-TextCell $ws.Cells.Item(38, 3) "3070"
Set-TextCell $ws.Cells.Item(38, 4) "Пошкодження поверхні контакту"
Set-TextCell $ws.Cells.Item(38, 5) "1207"

Set-TextCell $ws.Cells.Item(39, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(39, 2) "80001841"
Set-TextCell $ws.Cells.Item(39, 3) "3070"
Set-TextCell $ws.Cells.Item(39, 4) "Пошкодження поверхні контакту"
Set-TextCell $ws.Cells.Item(39, 5) "1207"

Set-TextCell $ws.Cells.Item(40, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(40, 2) "80001841"
Set-TextCell $ws.Cells.Item(40, 3) "3070"
Set-TextCell $ws.Cells.Item(40, 4) "Пошкодження поверхні контакту"
Set-TextCell $ws.Cells.Item(40, 5) "1208"

Set-TextCell $ws.Cells.Item(41, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(41, 2) "80001841"
Set-TextCell $ws.Cells.Item(41, 3) "3070"
Set-TextCell $ws.Cells.Item(41, 4) "Пошкодження поверхні контакту"
Set-TextCell $ws.Cells.Item(41, 5) "1209"

Set-TextCell $ws.Cells.Item(42, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(42, 2) "80001841"
Set-TextCell $ws.Cells.Item(42, 3) "3012"
Set-TextCell $ws.Cells.Item(42, 4) "Гострини на розрізі контакту"
Set-TextCell $ws.Cells.Item(42, 5) "1210"

Set-TextCell $ws.Cells.Item(43, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(43, 2) "80001841"
Set-TextCell $ws.Cells.Item(43, 3) "3012"
Set-TextCell $ws.Cells.Item(43, 4) "Гострини на розрізі контакту"
Set-TextCell $ws.Cells.Item(43, 5) "1211"

Set-TextCell $ws.Cells.Item(44, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(44, 2) "80001841"
Set-TextCell $ws.Cells.Item(44, 3) "3070"
Set-TextCell $ws.Cells.Item(44, 4) "Гострини на розрізі контакту"
Set-TextCell $ws.Cells.Item(44, 5) "1212"

Set-TextCell $ws.Cells.Item(45, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(45, 2) "80001841"
Set-TextCell $ws.Cells.Item(45, 3) "3070"
Set-TextCell $ws.Cells.Item(45, 4) "Гострини на розрізі контакту"
Set-TextCell $ws.Cells.Item(45, 5) "1220"

Set-TextCell $ws.Cells.Item(46, 1) "06/03/2018"
Set-TextCell $ws.Cells.Item(46, 2) "80001841"
Set-TextCell $ws.Cells.Item(46, 3) "3070"
Set-TextCell $ws.Cells.Item(46, 4) "Не вірна довжина проводу"
Set-TextCell $ws.Cells.Item(46, 5) "1230"

Set-TextCell $ws.Cells.Item(47, 1) "**"